# Applies the dades_GrupD.xlsx edit:
#  - Column I ("Salari mensual") values for rows 2..56 are re-scaled from
#    thousands (e.g. 2.611) to full units (e.g. 2611).
#  - Column T ("salario_mensual_actualizada") is recomputed as I * (1 + S)
#    where S is column S ("Increment"), matching the original relationship.
#  - Q17 ("Edad_actual") changes from 24 to 25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 56; $row++) {
    $iCell = $ws.Cells.Item($row, 9)   # column I
    $sCell = $ws.Cells.Item($row, 19)  # column S
    $tCell = $ws.Cells.Item($row, 20)  # column T

    $newSalary = [math]::Round([double]$iCell.Value2 * 1000, 6)
    $iCell.Value2 = $newSalary

    $increment = [double]$sCell.Value2
    $updatedSalary = [math]::Round($newSalary * (1 + $increment), 6)
    $tCell.Value2 = $updatedSalary
}

# Edad_actual for row 17 changes from 24 to 25
$ws.Cells.Item(17, 17).Value2 = 25
